# Update the "Buying Opportunity" table on the active sheet to reflect
# the latest scan results. Columns: A=index, B=Buying Opportunity,
# C=support Zone, D=long buildup, E=Short buildup, F=FII ENTERING.
# The table grows from 11 data rows (rows 2-12) to 17 data rows (rows 2-18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# New data set (rows 2..18 => index 0..16)
$data = @(
    @("NSE:APEX",         "NSE:ALANKIT",   "NSE:ATGL", "", "NSE:BANKBARODA"),
    @("NSE:BANKBARODA",   "NSE:BDL",       "",         "", "NSE:BANKINDIA"),
    @("NSE:BANKINDIA",    "NSE:CENTEXT",   "",         "", "NSE:BIOCON"),
    @("NSE:BIOCON",       "NSE:ENGINERSIN","",         "", "NSE:CONCOR"),
    @("NSE:BVCL",         "NSE:JYOTISTRUC","",         "", "NSE:LICHSGFIN"),
    @("NSE:DBCORP",       "NSE:LLOYDSENGG","",         "", "NSE:PATANJALI"),
    @("NSE:DPABHUSHAN",   "NSE:PFC",       "",         "", "NSE:PNB"),
    @("NSE:GARFIBRES",    "NSE:POONAWALLA","",         "", ""),
    @("NSE:HATHWAY",      "NSE:RTNINDIA",  "",         "", ""),
    @("NSE:KPIGREEN",     "NSE:RTNPOWER",  "",         "", ""),
    @("NSE:LPDC",         "",              "",         "", ""),
    @("NSE:MEDPLUS",      "",              "",         "", ""),
    @("NSE:NARMADA",      "",              "",         "", ""),
    @("NSE:ORIENTHOT",    "",              "",         "", ""),
    @("NSE:PATANJALI",    "",              "",         "", ""),
    @("NSE:PNB",          "",              "",         "", ""),
    @("NSE:PRAKASH",      "",              "",         "", "")
)

# Make sure the formatting (bold + border + centered/top aligned) used by
# the existing index column (A2) extends to the newly added rows (13-18)
# before we fill in the values.
$ws.Range("A2").Copy()
$ws.Range("A13:A18").PasteSpecial($xlPasteFormats)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
    $ws.Cells.Item($row, 6).Value = $data[$i][4]
}

Write-Host "Updated rows 2-18"
